# Actualización automática 2025-09-22 08:22:24
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": new sales figures for a couple of
# asesor/cliente rows (previously 0) plus the recalculated "X de 21"
# coverage counters on the totals row (row 23).
# ---------------------------------------------------------------------
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsVentasGrupo.Range("D5").Value = 1373.76
$wsVentasGrupo.Range("M5").Value = 3259.2
$wsVentasGrupo.Range("O5").Value = 547.43

$wsVentasGrupo.Range("H13").Value = 811.8
$wsVentasGrupo.Range("O13").Value = 547.43

$wsVentasGrupo.Range("D23").Value = "2 de 21"
$wsVentasGrupo.Range("H23").Value = "3 de 21"
$wsVentasGrupo.Range("M23").Value = "8 de 21"
$wsVentasGrupo.Range("O23").Value = "2 de 21"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL": "septiembre" column picks up the new sales for
# the same two rows, plus the updated column total.
# ---------------------------------------------------------------------
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsVentaMensual.Range("F5").Value = 5180.39
$wsVentaMensual.Range("F13").Value = 5833.17
$wsVentaMensual.Range("F23").Value = 44864.53

# ---------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL": widen the "VENTA" column slightly and
# refresh VENTA / POR CUMPLIR / CUMPLIMIENTO for the affected groups.
# ---------------------------------------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column E ("VENTA") width: XML `width` = ColumnWidth + 5/6, so to land
# on width="24" we set ColumnWidth to 24 - 5/6.
$wsCumplimiento.Columns.Item(5).ColumnWidth = 24 - 5/6

# Row 3 - 240X80 PORCELANATO
$wsCumplimiento.Range("D3").Value = 1848.96
$wsCumplimiento.Range("E3").Value = 3655.65890386263
$wsCumplimiento.Range("F3").Value = 0.3358924627284501

# Row 6 - INODOROS
$wsCumplimiento.Range("D6").Value = 2359.65
$wsCumplimiento.Range("E6").Value = 547.9336814602598
$wsCumplimiento.Range("F6").Value = 0.8115501593456894

# Row 12 - PORCELANATO
$wsCumplimiento.Range("D12").Value = 38098.83
$wsCumplimiento.Range("E12").Value = -1275.186907882904
$wsCumplimiento.Range("F12").Value = 1.034629569504922

# Row 14 - SAL SOLUBLE
$wsCumplimiento.Range("D14").Value = 1094.86
$wsCumplimiento.Range("E14").Value = -178.761404787216
$wsCumplimiento.Range("F14").Value = 1.195133368527538

# Row 15 - TOTAL
$wsCumplimiento.Range("D15").Value = 44864.53
$wsCumplimiento.Range("E15").Value = 10560.21316613377
$wsCumplimiento.Range("F15").Value = 0.8094675308737129
